$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-11-18 Saturday" "2023-11-19 Sunday"

Replace-Text "40×46=" "47×87="
Replace-Text "57×98=" "43×86="
Replace-Text "97×42=" "77×83="
Replace-Text "50×88=" "23×29="
Replace-Text "27×31=" "55×81="

Replace-Text "88×33=" "79×26="
Replace-Text "45×56=" "36×55="
Replace-Text "66×29=" "97×14="
Replace-Text "52×70=" "73×19="
Replace-Text "21×33=" "71×95="

Replace-Text "23×16=" "84×11="
Replace-Text "21×69=" "24×68="
Replace-Text "65×19=" "43×43="
Replace-Text "73×94=" "65×23="
Replace-Text "73×85=" "15×56="

Replace-Text "75×58=" "45×23="
Replace-Text "85×53=" "61×55="
Replace-Text "99×63=" "29×47="
Replace-Text "99×52=" "34×41="
Replace-Text "70×38=" "88×78="

Replace-Text "34×43=" "53×95="
Replace-Text "93×27=" "28×78="
Replace-Text "89×27=" "95×42="
Replace-Text "24×70=" "84×13="
Replace-Text "63×57=" "63×44="
